# Update "想去人数" (wanted-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 5,6,7,8,10,12,13 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 16
$wsExhibit.Range("F6").Value = 3197
$wsExhibit.Range("F7").Value = 2096
$wsExhibit.Range("F8").Value = 401
$wsExhibit.Range("F10").Value = 1188
$wsExhibit.Range("F12").Value = 1102
$wsExhibit.Range("F13").Value = 90

# Sheet "全部类型": rows 5,6,7,8,11,13,14 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 16
$wsAll.Range("F6").Value = 3197
$wsAll.Range("F7").Value = 2096
$wsAll.Range("F8").Value = 401
$wsAll.Range("F11").Value = 1188
$wsAll.Range("F13").Value = 1102
$wsAll.Range("F14").Value = 90
